$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-14 Sunday", "2025-12-15 Monday"),
    @("595×3=1785", "137×2=274"),
    @("954×5=4770", "538×5=2690"),
    @("525×7=3675", "155×7=1085"),
    @("722×5=3610", "906×3=2718"),
    @("249×4=996", "205×2=410"),
    @("391×5=1955", "226×8=1808"),
    @("891×9=8019", "111×7=777"),
    @("913×7=6391", "152×4=608"),
    @("250×6=1500", "718×9=6462"),
    @("362×8=2896", "561×5=2805"),
    @("731×7=5117", "166×8=1328"),
    @("631×7=4417", "882×6=5292"),
    @("142×4=568", "306×5=1530"),
    @("980×9=8820", "556×7=3892"),
    @("331×3=993", "476×7=3332"),
    @("862×8=6896", "365×7=2555"),
    @("474×6=2844", "964×5=4820"),
    @("122×4=488", "559×4=2236"),
    @("862×9=7758", "541×3=1623"),
    @("620×2=1240", "279×5=1395"),
    @("867×6=5202", "994×6=5964"),
    @("269×8=2152", "640×9=5760"),
    @("562×3=1686", "550×9=4950"),
    @("672×4=2688", "555×2=1110"),
    @("981×9=8829", "407×9=3663")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
